# Append a new paragraph at the end of the document body, after the
# paragraph ending "... which is not a good use of time."
$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Insert a new paragraph mark after the current last paragraph; the new
# paragraph inherits the preceding paragraph's formatting (Arial font,
# justified alignment), matching the target OOXML's <w:pPr>/<w:rPr>.
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Text = "Most interestingly, this review revealed other small things that only individual group members noticed meaning they could share them with the group."
